$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 (Cx.m) - establishes shared strings "angle (rads)" and "3x3 (float)" first
$ws.Range("B16").Value = "angle (rads)"
$ws.Range("C16").Value = "3x3 (float)"
$ws.Range("E16").Value = "Tom"
$ws.Range("F16").Value = 1
$ws.Range("F16").NumberFormat = $ws.Range("F2").NumberFormat

# Row 17 (Cy.m)
$ws.Range("B17").Value = "angle (rads)"
$ws.Range("C17").Value = "3x3 (float)"
$ws.Range("E17").Value = "Tom"
$ws.Range("F17").Value = 1
$ws.Range("F17").NumberFormat = $ws.Range("F2").NumberFormat

# Row 18 (Cz.m)
$ws.Range("B18").Value = "angle (rads)"
$ws.Range("C18").Value = "3x3 (float)"
$ws.Range("E18").Value = "Tom"
$ws.Range("F18").Value = 1
$ws.Range("F18").NumberFormat = $ws.Range("F2").NumberFormat

# Row 11 (StateRates.m)
$ws.Range("E11").Value = "Tom"

# Row 3 (Trim.m)
$ws.Range("E3").Value = "Tom"
$ws.Range("B3").Value = "state vector"
$ws.Range("C3").Value = "trim inputs"
$ws.Range("F3").Value = 0
$ws.Range("F3").NumberFormat = $ws.Range("F2").NumberFormat

# Update selection to match new active cell
$ws.Range("E13").Select()
